# Update the cryptos list worksheet with the latest scraped values.
# All of these columns are stored as literal text (inline strings) in the
# original workbook, even when a value "looks like" a number (e.g. prices
# such as "0.8900" or "301.50" where trailing zeros matter, or big numbers
# like "26.825.74" that use '.' as a thousands separator). Excel's COM
# automation auto-coerces a plain numeric-looking string assigned via
# Range.Value into a real number (dropping the formatting/trailing
# zeros), so for any new value that parses as a plain number we force the
# cell to Text format first, assign the literal string, then clear the
# formatting back off (which also resets the style index) so the cell
# keeps its original "General" appearance while retaining the exact text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = '@'
    $r.Value = $val
    $r.ClearFormats()
}

$ws.Range('D2').Value = '26.825.74'
$ws.Range('E2').Value = '  -1.49%  '
$ws.Range('D3').Value = '1.874.14'
$ws.Range('E3').Value = '  -1.72%  '
Set-TextValue 'D4' '1.002'
$ws.Range('E4').Value = '  -0.05%  '
Set-TextValue 'D5' '301.50'
$ws.Range('E5').Value = '  -2.03%  '
$ws.Range('E6').Value = '  -0.02%  '
Set-TextValue 'D7' '0.5369'
$ws.Range('E7').Value = '  +2.23%  '
Set-TextValue 'D8' '0.3745'
$ws.Range('E8').Value = '  -2.01%  '
Set-TextValue 'D9' '0.07197'
$ws.Range('E9').Value = '  -1.53%  '
Set-TextValue 'D10' '21.59'
$ws.Range('E10').Value = '  +0.08%  '
Set-TextValue 'D11' '0.8900'
$ws.Range('E11').Value = '  -1.75%  '
Set-TextValue 'D12' '0.08170'
$ws.Range('E12').Value = '  +0.99%  '
$ws.Range('D13').Value = '1.878.08'
$ws.Range('E13').Value = '  +6.47%  '
Set-TextValue 'D14' '93.45'
$ws.Range('E14').Value = '  -2.69%  '
Set-TextValue 'D15' '5.322'
$ws.Range('E15').Value = '  -0.91%  '
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('E17').Value = '  +0.74%  '
Set-TextValue 'D18' '0.000008537'
$ws.Range('E18').Value = '  -1.58%  '
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').Value = '26.864.28'
$ws.Range('E20').Value = '  -1.49%  '
Set-TextValue 'D21' '4.989'
$ws.Range('E21').Value = '  -2.58%  '
$ws.Range('E22').Value = '  -1.79%  '
Set-TextValue 'D23' '6.404'
$ws.Range('E23').Value = '  -1.33%  '
Set-TextValue 'D24' '2.291'
$ws.Range('E24').Value = '  -2.44%  '
Set-TextValue 'D25' '146.41'
$ws.Range('E25').Value = '  -2.46%  '
Set-TextValue 'D26' '18.11'
$ws.Range('E26').Value = '  -0.80%  '
Set-TextValue 'D27' '1.732'
$ws.Range('E27').Value = '  -0.53%  '
Set-TextValue 'D28' '114.06'
$ws.Range('E28').Value = '  -2.37%  '
Set-TextValue 'D29' '4.722'
$ws.Range('E29').Value = '  -2.79%  '
Set-TextValue 'D30' '4.614'
$ws.Range('E30').Value = '  -5.48%  '
Set-TextValue 'D31' '0.09140'
$ws.Range('E31').Value = '  -0.91%  '
Set-TextValue 'D32' '0.8141'
$ws.Range('E32').Value = '  -1.39%  '
Set-TextValue 'D33' '0.05024'
Set-TextValue 'D34' '1.175'
$ws.Range('E34').Value = '  -4.63%  '
Set-TextValue 'D35' '2.952'
$ws.Range('E35').Value = '  -1.43%  '
Set-TextValue 'D36' '0.6036'
$ws.Range('E36').Value = '  +4.99%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D37' '3.220'
$ws.Range('E37').Value = '  -4.28%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D38' '2.628'
$ws.Range('E38').Value = '  -3.93%  '
Set-TextValue 'D39' '0.01956'
$ws.Range('E39').Value = '  -2.45%  '
Set-TextValue 'D40' '1.071'
$ws.Range('E40').Value = '  -1.15%  '
Set-TextValue 'D41' '6.638'
$ws.Range('E41').Value = '  +0.23%  '
Set-TextValue 'D42' '8.903'
$ws.Range('E42').Value = '  -1.73%  '
Set-TextValue 'D43' '115.26'
$ws.Range('E43').Value = '  -1.75%  '
Set-TextValue 'D44' '0.5116'
$ws.Range('E44').Value = '  +3.62%  '
Set-TextValue 'D45' '0.1497'
$ws.Range('E45').Value = '  -1.79%  '
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D47' '1.645'
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D48' '9.971'
$ws.Range('E48').Value = '  -2.50%  '
Set-TextValue 'D49' '37.67'
$ws.Range('E49').Value = '  -2.58%  '
Set-TextValue 'D50' '0.06083'
$ws.Range('E50').Value = '  +1.87%  '
Set-TextValue 'D51' '62.23'
$ws.Range('E51').Value = '  -3.12%  '
